# Qual_6_data.xlsx update:
# - Fill in the previously-missing TIME_PERIOD (column L) value of 2019
#   for the batch of rows (192-382) that were added/sourced in 2019 but
#   had never had their TIME_PERIOD back-filled.
# - The indicator/series range now spans the full A:O columns (instead of
#   the old E:N), so the sheet's AutoFilter and the hidden
#   "_xlnm._FilterDatabase" defined name are widened to match.
# - Re-select the cell the editor ended up on after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Back-fill TIME_PERIOD (column L = 12) with 2019 for rows 192-382 ---
for ($r = 192; $r -le 382; $r++) {
    $ws.Cells.Item($r, 12).Value = 2019
}

# --- 2. Widen the hidden _xlnm._FilterDatabase defined name to A1:O382 ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='6 - Presence of a gender quota '!`$A`$1:`$O`$382"
    }
}

# --- 3. Widen the AutoFilter range to A1:O382 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:O382").AutoFilter()

# --- 4. Leave the selection on the last touched cell ---
$ws.Range("L378").Select()
